# Updated main GSC export data: the oldest day (2025-10-19, row 2 of the
# "Chart" sheet) has rolled off the reporting window. Delete that entire
# row so every subsequent row shifts up by one; Excel automatically
# removes the now-unused "2025-10-19" shared string and keeps the
# "Table" sheet's header references consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
